$d = $word.ActiveDocument

# --- Paragraph 1: "Fernando Viana Junior" ---
$p1 = $d.Paragraphs.Item(1)

# Remove center alignment (jc w:val="center") -> left/default alignment
$p1.Alignment = 0

# Append "!" as its own run (split it off from the existing run by
# briefly toggling bold so Word has to materialize a new <w:r>, then
# toggling it back off so formatting matches the surrounding text).
$insertPoint = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$insertPoint.InsertAfter("!")
$bangRange = $d.Range($p1.Range.End - 2, $p1.Range.End - 1)
$bangRange.Font.Bold = 1
$bangRange.Font.Bold = 0

# --- Paragraph 2: "1460681823018" ---
$p2 = $d.Paragraphs.Item(2)

# Remove center alignment (jc w:val="center") -> left/default alignment
$p2.Alignment = 0

# The paragraph currently starts with a line break followed by the
# number: <w:br/><w:t>1460681823018</w:t>. Remove that leading break.
$leadingBreak = $d.Range($p2.Range.Start, $p2.Range.Start + 1)
$leadingBreak.Delete()

# Insert "." right before the paragraph mark (merges into the number's
# run, inheriting its formatting/rFonts).
$paraMark = $d.Range($p2.Range.End - 1, $p2.Range.End)
$paraMark.InsertBefore(".")

# Insert a line break right before the paragraph mark too, so the
# paragraph ends with ...018.<br/>
$paraMark2 = $d.Range($p2.Range.End - 1, $p2.Range.End)
$paraMark2.InsertBefore([char]11)

# Split "." and the trailing break into their own runs (toggle bold
# on/off so formatting ends up identical to their neighbors but Word
# still emits separate <w:r> elements).
$breakRun = $d.Range($p2.Range.End - 2, $p2.Range.End - 1)
$breakRun.Font.Bold = 1
$breakRun.Font.Bold = 0

$dotRun = $d.Range($p2.Range.End - 3, $p2.Range.End - 2)
$dotRun.Font.Bold = 1
$dotRun.Font.Bold = 0
